$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-wise new values for columns E, F, G, H (common to all data rows 2-6)
$rows = @(2, 3, 4, 5, 6)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = 2
    $ws.Range("F$r").Value = 1
    $ws.Range("G$r").Value = 0.307034
    $ws.Range("H$r").Value = 0.6140680000000001
}

# Row 2 specific updates
$ws.Range("M2").Value = 14.37161333333333
$ws.Range("N2").Value = 43.11484
$ws.Range("O2").Value = 0.4561705932627708
$ws.Range("P2").Value = 0.5019766122855294
$ws.Range("Q2").Value = 4.412573928186667
$ws.Range("R2").Value = 26.47544356912
$ws.Range("S2").Value = 0.4561705932627708
$ws.Range("T2").Value = 0.5019766122855294

# Row 3 specific updates
$ws.Range("O3").Value = 0.01117178254830525
$ws.Range("P3").Value = 0.01229358849433434
$ws.Range("Q3").Value = 0.108065528844
$ws.Range("R3").Value = 0.648393173064
$ws.Range("S3").Value = 0.01117178254830525
$ws.Range("T3").Value = 0.01229358849433434

# Row 4 specific updates
$ws.Range("M4").Value = 2.798424666666667
$ws.Range("N4").Value = 8.395274000000001
$ws.Range("O4").Value = 0.08882503382091908
$ws.Range("P4").Value = 0.09774433122629669
$ws.Range("Q4").Value = 0.8592115191053336
$ws.Range("R4").Value = 5.155269114632
$ws.Range("S4").Value = 0.08882503382091908
$ws.Range("T4").Value = 0.09774433122629669

# Row 5 specific updates
$ws.Range("M5").Value = 8.624592
$ws.Range("N5").Value = 17.249184
$ws.Range("O5").Value = 0.2737539034788959
$ws.Range("P5").Value = 0.2008284606648142
$ws.Range("Q5").Value = 2.648042980128
$ws.Range("R5").Value = 10.592171920512
$ws.Range("S5").Value = 0.2737539034788959
$ws.Range("T5").Value = 0.2008284606648142

# Row 6 specific updates
$ws.Range("M6").Value = 5.358313666666667
$ws.Range("N6").Value = 16.074941
$ws.Range("O6").Value = 0.1700786868891091
$ws.Range("P6").Value = 0.1871570073290255
$ws.Range("Q6").Value = 1.645184478331333
$ws.Range("R6").Value = 9.871106869988001
$ws.Range("S6").Value = 0.1700786868891091
$ws.Range("T6").Value = 0.1871570073290255
